$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 5 holds the data that should become the new row 6
# (duplicated down), while row 5 itself gets updated with new values.
# Copy row 5 -> row 6 cell by cell (values only) before changing row 5.

for ($col = 1; $col -le 18; $col++) {
    $srcCell = $ws.Cells.Item(5, $col)
    $dstCell = $ws.Cells.Item(6, $col)
    $dstCell.Value = $srcCell.Value()
}

# Match row 5's style on row 6 (covers the date number format on column D).
$ws.Range("A5:R5").Copy()
$ws.Range("A6:R6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Now update row 5 in place with the new values from the diff.
$ws.Cells.Item(5, 4).Value = 44755    # D5 Fecha
$ws.Cells.Item(5, 10).Value = 50      # J5 Volumen
$ws.Cells.Item(5, 11).Value = 20000   # K5 Precio minimo
$ws.Cells.Item(5, 12).Value = 20000   # L5 Precio maximo
$ws.Cells.Item(5, 13).Value = 20000   # M5 Precio promedio ponderado
$ws.Cells.Item(5, 16).Value = 1333    # P5 Precio $/Kg
